$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.315.50'
$ws.Range("E2").Value = '  +0.29%  '

$ws.Range("D3").Value = '1.876.86'
$ws.Range("E3").Value = '  +0.87%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = "'0.7120"
$ws.Range("E5").Value = '  -0.54%  '

$ws.Range("D6").Value = "'242.32"
$ws.Range("E6").Value = '  +0.65%  '

$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("D8").Value = "'0.3106"
$ws.Range("E8").Value = '  +0.72%  '

$ws.Range("D9").Value = "'0.07770"
$ws.Range("E9").Value = '  +0.24%  '

$ws.Range("D10").Value = "'24.97"
$ws.Range("E10").Value = '  -0.92%  '

$ws.Range("E11").Value = '  +2.47%  '

$ws.Range("D12").Value = '1.885.56'
$ws.Range("E12").Value = '  +0.85%  '

$ws.Range("E13").Value = '  -0.46%  '

$ws.Range("E14").Value = '  -0.97%  '

$ws.Range("D15").Value = "'91.47"
$ws.Range("E15").Value = '  +1.29%  '

$ws.Range("D16").Value = '29.315.88'
$ws.Range("E16").Value = '  +0.39%  '

$ws.Range("D17").Value = "'0.000008264"
$ws.Range("E17").Value = '  +5.81%  '

$ws.Range("E18").Value = '  +2.38%  '

$ws.Range("D19").Value = "'242.45"
$ws.Range("E19").Value = '  -0.87%  '

$ws.Range("D20").Value = '2.134.24'
$ws.Range("E20").Value = '  +1.73%  '

$ws.Range("D21").Value = "'13.26"
$ws.Range("E21").Value = '  +0.74%  '

$ws.Range("E22").Value = '  -0.04%  '

$ws.Range("D23").Value = "'7.836"
$ws.Range("E23").Value = '  -1.57%  '

$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = '  +0.04%  '

$ws.Range("D25").Value = "'0.1618"
$ws.Range("E25").Value = '  +1.46%  '

$ws.Range("D26").Value = "'162.64"
$ws.Range("E26").Value = '  -0.11%  '

$ws.Range("D27").Value = "'9.029"
$ws.Range("E27").Value = '  +0.95%  '

$ws.Range("D29").Value = "'1.511"
$ws.Range("E29").Value = '  +1.00%  '

$ws.Range("D30").Value = "'4.410"
$ws.Range("E30").Value = '  +0.11%  '

$ws.Range("E31").Value = '  +4.05%  '

$ws.Range("E32").Value = '  -2.80%  '

$ws.Range("D33").Value = "'0.05239"
$ws.Range("E33").Value = '  +0.71%  '

$ws.Range("E34").Value = '  +1.14%  '

$ws.Range("D35").Value = "'1.179"
$ws.Range("E35").Value = '  +0.47%  '

$ws.Range("D36").Value = "'0.7412"
$ws.Range("E36").Value = '  +1.81%  '

$ws.Range("D37").Value = "'2.686"
$ws.Range("E37").Value = '  +0.23%  '

$ws.Range("D38").Value = "'0.01864"
$ws.Range("E38").Value = '  +0.46%  '

$ws.Range("D39").Value = "'2.725"
$ws.Range("E39").Value = '  +1.49%  '

$ws.Range("D40").Value = '1.173.65'
$ws.Range("E40").Value = '  +1.88%  '

$ws.Range("D41").Value = "'6.389"
$ws.Range("E41").Value = '  +4.74%  '

$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = "'72.97"
$ws.Range("E42").Value = '  +0.74%  '

$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = "'0.8878"
$ws.Range("E43").Value = '  -2.20%  '

$ws.Range("D44").Value = "'106.51"
$ws.Range("E44").Value = '  +4.66%  '

$ws.Range("E45").Value = '  -0.01%  '

$ws.Range("D46").Value = '2.029.29'
$ws.Range("E46").Value = '  +1.26%  '

$ws.Range("E47").Value = '  +2.62%  '

$ws.Range("D48").Value = "'0.5204"
$ws.Range("E48").Value = '  -0.60%  '

$ws.Range("E49").Value = '  +4.82%  '

$ws.Range("D50").Value = "'9.410"

$ws.Range("D51").Value = "'0.4314"
$ws.Range("E51").Value = '  +1.05%  '
